$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/12/2024  Through  8/18/2024"

# --- Cells converting from "n/a" placeholder text to real numbers ---
$ws.Range("C16").Value = 1
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("D16").Value = 2
$ws.Range("D16").NumberFormat = '#,##0'
$ws.Range("E16").Value = -50
$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C18").Value = 1
$ws.Range("C18").NumberFormat = '#,##0'
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("F22").Value = 1
$ws.Range("F22").NumberFormat = '#,##0'

# --- Cells converting from real numbers back to "n/a" placeholder text ---
# (copy format+value from an existing placeholder cell with identical style)
$ws.Range("C14").Copy($ws.Range("G33"))   # "0" placeholder (s=14, shared string 20)
$ws.Range("E14").Copy($ws.Range("H33"))   # "***.*" placeholder (s=14, shared string 21)

# --- Plain numeric value updates ---
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 42
$ws.Range("J16").Value = 46
$ws.Range("K16").Value = -8.695652173913
$ws.Range("L16").Value = 35.483870967741
$ws.Range("M16").Value = -33.333333333333
$ws.Range("N16").Value = -88
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 9.090909090909
$ws.Range("I17").Value = 77
$ws.Range("J17").Value = 92
$ws.Range("K17").Value = -16.304347826087
$ws.Range("L17").Value = 11.594202898550
$ws.Range("M17").Value = 16.666666666666
$ws.Range("N17").Value = -57.458563535911
$ws.Range("E18").Value = -50
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -57.142857142857
$ws.Range("I18").Value = 50
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = -24.242424242424
$ws.Range("L18").Value = -1.960784313725
$ws.Range("M18").Value = -67.105263157894
$ws.Range("N18").Value = -92.603550295858
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -72.727272727272
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -41.176470588235
$ws.Range("I19").Value = 240
$ws.Range("J19").Value = 278
$ws.Range("K19").Value = -13.669064748201
$ws.Range("L19").Value = -29.411764705882
$ws.Range("M19").Value = 20.603015075376
$ws.Range("N19").Value = -17.808219178082
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 131
$ws.Range("J20").Value = 80
$ws.Range("K20").Value = 63.75
$ws.Range("L20").Value = 77.027027027027
$ws.Range("M20").Value = 25.961538461538
$ws.Range("N20").Value = -89.486356340288
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -36.363636363636
$ws.Range("F21").Value = 56
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = -22.222222222222
$ws.Range("I21").Value = 545
$ws.Range("J21").Value = 570
$ws.Range("K21").Value = -4.385964912280
$ws.Range("L21").Value = -5.052264808362
$ws.Range("M21").Value = -7.783417935702
$ws.Range("N21").Value = -80.282199710564
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = 16.666666666666
$ws.Range("L22").Value = -30
$ws.Range("M22").Value = 16.666666666666
$ws.Range("C24").Value = 45
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 152
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = 44.761904761904
$ws.Range("I24").Value = 914
$ws.Range("J24").Value = 958
$ws.Range("K24").Value = -4.592901878914
$ws.Range("L24").Value = -21.070811744386
$ws.Range("M24").Value = 25.722145804676
$ws.Range("C25").Value = 28
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 64.705882352941
$ws.Range("F25").Value = 90
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = 83.673469387755
$ws.Range("I25").Value = 582
$ws.Range("J25").Value = 479
$ws.Range("K25").Value = 21.503131524008
$ws.Range("L25").Value = -14.537444933920
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 114.285714285714
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = 31.818181818181
$ws.Range("I26").Value = 228
$ws.Range("J26").Value = 235
$ws.Range("K26").Value = -2.978723404255
$ws.Range("L26").Value = 10.144927536231
$ws.Range("M26").Value = -0.869565217391
$ws.Range("C28").Value = 5
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 600
$ws.Range("I28").Value = 34
$ws.Range("K28").Value = 88.888888888888
$ws.Range("L28").Value = 9.677419354838
